$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date column (B) for rows 2-7 with six new timestamp values,
# replacing the previous "Wed Feb 15 12:4x:xx EST 2023" strings with the
# new "Thu Feb 16 12:1x:xx EST 2023" strings (new shared strings get
# appended to the shared string table).
$ws.Range("B2").Value = "Thu Feb 16 12:10:45 EST 2023"
$ws.Range("B3").Value = "Thu Feb 16 12:10:55 EST 2023"
$ws.Range("B4").Value = "Thu Feb 16 12:11:05 EST 2023"
$ws.Range("B5").Value = "Thu Feb 16 12:11:14 EST 2023"
$ws.Range("B6").Value = "Thu Feb 16 12:11:24 EST 2023"
$ws.Range("B7").Value = "Thu Feb 16 12:11:34 EST 2023"
